$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 10.99600299506004
    "C2" = 4.563921367383737
    "D2" = 6.005904277717514
    "E2" = 16.32460510360481
    "G2" = 52.02437607013296
    "H2" = 19.15561665238971
    "K2" = 10.34339186840125
    "B3" = 10.77918560200353
    "C3" = 4.30842448648003
    "D3" = 5.891730649253343
    "E3" = 15.40753757723155
    "G3" = 50.92497259606557
    "H3" = 19.02101221183007
    "K3" = 10.21577666198213
    "B4" = 10.64852265580841
    "C4" = 4.142733173993973
    "D4" = 5.822420649956568
    "E4" = 14.82120478240416
    "G4" = 50.24663311047358
    "H4" = 18.94047329172316
    "K4" = 10.14093459050683
    "B5" = 10.59598624054826
    "C5" = 4.073010223510128
    "D5" = 5.794417386109912
    "E5" = 14.57669819818201
    "G5" = 49.96970954448911
    "H5" = 18.90820120737083
    "K5" = 10.11136151351496
    "B6" = 10.58730810696556
    "C6" = 4.061300263019662
    "D6" = 5.789783247781333
    "E6" = 14.53577066074851
    "G6" = 49.9237066707941
    "H6" = 18.90287606421656
    "K6" = 10.10650801744212
    "B7" = 10.64781113853636
    "C7" = 4.141801761294735
    "D7" = 5.822041956696085
    "E7" = 14.81792943615548
    "G7" = 50.24289997676073
    "H7" = 18.94003581680762
    "K7" = 10.14053195589641
    "B8" = 10.92078719838816
    "C8" = 4.477662909338251
    "D8" = 5.966397584906487
    "E8" = 16.0133750789141
    "G8" = 51.64618554454717
    "H8" = 19.1087784147281
    "K8" = 10.2986866557343
    "B9" = 11.47168561601716
    "C9" = 5.065915032788715
    "D9" = 6.253953647895713
    "E9" = 18.19302772853116
    "G9" = 54.35757779880566
    "H9" = 19.45558566185493
    "K9" = 10.63476283632033
    "B10" = 11.88073148325075
    "C10" = 5.454862089168173
    "D10" = 6.465566678072547
    "E10" = 19.82059442274585
    "G10" = 56.30674924880827
    "H10" = 19.71891790440569
    "K10" = 10.8947716826444
    "B11" = 12.06669782233134
    "C11" = 5.622392039117742
    "D11" = 6.561428865926168
    "E11" = 20.51981398594281
    "G11" = 57.18059117984208
    "H11" = 19.84030492199061
    "K11" = 11.01530473896904
    "B12" = 12.13701966981779
    "C12" = 5.684479048460781
    "D12" = 6.59763370917563
    "E12" = 20.77870875923922
    "G12" = 57.50936876830872
    "H12" = 19.88647736778803
    "K12" = 11.06122166482949
    "B13" = 12.12188031615676
    "C13" = 5.671167668151226
    "D13" = 6.589841211792641
    "E13" = 20.72321228885398
    "G13" = 57.43865952161011
    "H13" = 19.87652453290089
    "K13" = 11.0513212590451
    "B14" = 12.07248570242266
    "C14" = 5.627527088492402
    "D14" = 6.564409598637668
    "E14" = 20.54123114611755
    "G14" = 57.20768408806168
    "H14" = 19.84409956907605
    "K14" = 11.01907715433594
    "B15" = 12.04221469886689
    "C15" = 5.600619784816052
    "D15" = 6.548818382885798
    "E15" = 20.42899715684861
    "G15" = 57.06591990256474
    "H15" = 19.8242644362501
    "K15" = 10.99936086119466
    "B16" = 11.86856867099494
    "C16" = 5.443724289835024
    "D16" = 6.459290456891405
    "E16" = 19.77407339989368
    "G16" = 56.24935878259636
    "H16" = 19.71101507714534
    "K16" = 10.88693559640306
    "B17" = 11.76195326366401
    "C17" = 5.345064442324281
    "D17" = 6.404237239367319
    "E17" = 19.36177992014701
    "G17" = 55.74492997076423
    "H17" = 19.64193230759174
    "K17" = 10.81850759325214
    "B18" = 11.70062523524595
    "C18" = 5.287432691546303
    "D18" = 6.372536459917042
    "E18" = 19.12076251091425
    "G18" = 55.45360240605704
    "H18" = 19.60234938384579
    "K18" = 10.77936475319018
    "B19" = 11.67986226322789
    "C19" = 5.267767578962614
    "D19" = 6.361798200777918
    "E19" = 19.03849117670236
    "G19" = 55.35476774380164
    "H19" = 19.588974047095
    "K19" = 10.76615015225729
    "B20" = 11.77330386176831
    "C20" = 5.35565857725477
    "D20" = 6.410101719024267
    "E20" = 19.40607027156156
    "G20" = 55.79875288193972
    "H20" = 19.64927076057866
    "K20" = 10.82576997170259
    "B21" = 12.08699744543571
    "C21" = 5.640382101339126
    "D21" = 6.571882383646852
    "E21" = 20.59484284057877
    "G21" = 57.27558705393569
    "H21" = 19.85361817181406
    "K21" = 11.02854100147479
    "B22" = 12.29139828375079
    "C22" = 5.818580029053991
    "D22" = 6.677038935911522
    "E22" = 21.33749141913158
    "G22" = 58.22826607900848
    "H22" = 19.98835942651901
    "K22" = 11.16263580324891
    "B23" = 12.18238858576808
    "C23" = 5.724193810338394
    "D23" = 6.620979813017523
    "E23" = 20.94425110825627
    "G23" = 57.72103669266816
    "H23" = 19.91634466468075
    "K23" = 11.09093953946714
    "B24" = 11.76817235398795
    "C24" = 5.350871806093598
    "D24" = 6.40745054045638
    "E24" = 19.3860590174495
    "G24" = 55.77442364852945
    "H24" = 19.64595262768627
    "K24" = 10.8224860318737
    "B25" = 11.32152948090703
    "C25" = 4.914356355963855
    "D25" = 6.175933139834844
    "E25" = 17.60380151682253
    "G25" = 53.63033987670057
    "H25" = 19.36017935946257
    "K25" = 10.54134518689225
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
